$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns keep their original text formatting
# by forcing a text number format before assigning values that could otherwise
# be auto-converted by Excel into numeric values (e.g. "1.00" -> 1).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.609.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.915.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.36"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +8.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.40"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.614"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.34%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.723"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.73%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000335"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.23"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.523.75"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.28"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.908.74"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.14"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.60%  "

$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.22"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +7.05%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.135"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.75"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.503.98"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "429.20"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.36"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.56%  "

$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.25"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.95%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.50"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.02"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +8.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.49"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.60"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.42"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "695.82"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.21"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.126"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.82"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "67.32"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +10.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.439"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +9.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.94"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.15"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0845"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.12%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0481"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.11"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.72%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.78"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -9.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.05"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.73%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.141"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.33"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.80%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.03"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.57%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0353"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.95%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.67"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.48%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.724.34"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +11.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.29"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.98%  "
